# Update "想去人数" (interested-count) figures to the freshly generated data.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8902
$ws1.Range("F4").Value = 238
$ws1.Range("F6").Value = 1488
$ws1.Range("F7").Value = 1409
$ws1.Range("F10").Value = 319

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 16

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8902
$ws4.Range("F4").Value = 238
$ws4.Range("F6").Value = 1488
$ws4.Range("F7").Value = 1409
$ws4.Range("F9").Value = 16
$ws4.Range("F11").Value = 319
